$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at position 236, shifting the existing rows
# (old rows 236..343) down to 238..345.
$ws.Rows("236:237").Insert()

# The two newly inserted (blank) rows receive the data that used to
# live in rows 234 and 235 before this edit (row 234 -> new row 236,
# row 235 -> new row 237).
$ws.Range("A236").Value2 = 5
$ws.Range("B236").Value2 = "Macroferia Regional de Talca"
$ws.Range("C236").Value2 = "Maule"
$ws.Range("D236").Value2 = 44389
$ws.Range("E236").Value2 = 7
$ws.Range("F236").Value2 = 100112020
$ws.Range("G236").Value2 = "Tomate"
$ws.Range("H236").Value2 = "Larga vida"
$ws.Range("I236").Value2 = "Primera"
$ws.Range("J236").Value2 = 1500
$ws.Range("K236").Value2 = 12000
$ws.Range("L236").Value2 = 12000
$ws.Range("M236").Value2 = 12000
$ws.Range("N236").Value2 = "`$/bandeja 18 kilos"
$ws.Range("O236").Value2 = "Región de Arica y Parinacota"
$ws.Range("P236").Value2 = 667
$ws.Range("Q236").Value2 = 18
$ws.Range("R236").Value2 = "Hortaliza"

$ws.Range("A237").Value2 = 5
$ws.Range("B237").Value2 = "Macroferia Regional de Talca"
$ws.Range("C237").Value2 = "Maule"
$ws.Range("D237").Value2 = 44389
$ws.Range("E237").Value2 = 7
$ws.Range("F237").Value2 = 100112020
$ws.Range("G237").Value2 = "Tomate"
$ws.Range("H237").Value2 = "Larga vida"
$ws.Range("I237").Value2 = "Primera"
$ws.Range("J237").Value2 = 1500
$ws.Range("K237").Value2 = 6000
$ws.Range("L237").Value2 = 6000
$ws.Range("M237").Value2 = 6000
$ws.Range("N237").Value2 = "`$/caja 10 kilos"
$ws.Range("O237").Value2 = "Región de Arica y Parinacota"
$ws.Range("P237").Value2 = 600
$ws.Range("Q237").Value2 = 10
$ws.Range("R237").Value2 = "Hortaliza"

# Rows 234 and 235 now hold new data for the 2021-09-20 (44466) date.
$ws.Range("D234").Value2 = 44466
$ws.Range("K234").Value2 = 14000
$ws.Range("L234").Value2 = 14000
$ws.Range("M234").Value2 = 14000
$ws.Range("P234").Value2 = 778

$ws.Range("D235").Value2 = 44466
